$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.930.52'
$ws.Range('E2').Value = '  -3.37%  '
$ws.Range('D3').Value = '3.019.67'
$ws.Range('E3').Value = '  -4.12%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -6.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '3.016.28'
$ws.Range('E8').Value = '  -3.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.148'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.01'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.25%  '
$ws.Range('E12').Value = '  -4.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000219'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.13'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -8.59%  '
$ws.Range('D15').Value = '3.510.09'
$ws.Range('E15').Value = '  -4.08%  '
$ws.Range('D16').Value = '61.861.16'
$ws.Range('E16').Value = '  -3.72%  '
$ws.Range('E17').Value = '  -2.62%  '
$ws.Range('D18').Value = '3.015.68'
$ws.Range('E18').Value = '  -4.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '469.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.679'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -9.66%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '25.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -13.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.08'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '55.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -10.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.76'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '460.38'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -16.29%  '
$ws.Range('D38').Value = '3.026.65'
$ws.Range('E38').Value = '  -4.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0381'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0771'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.111'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.86'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.46'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.242'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.50%  '
$ws.Range('B46').Value = 'PEPE'
$ws.Range('C46').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D46').Value = '0.0₃0520'
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '117.54'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -10.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.106'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.30'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.68%  '
